$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: new data row, values stored as text (matches source inlineStr "5000.0"/"4000.0"/"20")
$row2 = $ws.Range("A2:C2")
$row2.NumberFormat = "@"

$ws.Range("A2").Value = "5000.0"
$ws.Range("B2").Value = "4000.0"
$ws.Range("C2").Value = "20"

# Restore default ("Normal") cell style so the new row matches the un-styled
# header-less data cells from the source (no custom formatting carried over).
$row2.Style = "Normal"
